$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 1.037825649916365
$block1[0,2] = 1.044782224089829
$block1[0,3] = 1.036548577346475
$block1[0,4] = 1.052173965504242
$block1[1,0] = 1.02
$block1[1,1] = 1.039011970126322
$block1[1,2] = 1.04587962936566
$block1[1,3] = 1.037563782887696
$block1[1,4] = 1.053434890126086
$block1[2,0] = 1.02
$block1[2,1] = 1.039779289200543
$block1[2,2] = 1.046589759672108
$block1[2,3] = 1.038220752955098
$block1[2,4] = 1.054251126428484
$block1[3,0] = 1.02
$block1[3,1] = 1.040101797786183
$block1[3,2] = 1.046888308035674
$block1[3,3] = 1.038496959280671
$block1[3,4] = 1.054594353520268
$block1[4,0] = 1.02
$block1[4,1] = 1.040155944204103
$block1[4,2] = 1.046938436211519
$block1[4,3] = 1.038543336504213
$block1[4,4] = 1.054651987671746
$block1[5,0] = 1.02
$block1[5,1] = 1.039783598858197
$block1[5,2] = 1.046593748851016
$block1[5,3] = 1.038224443573325
$block1[5,4] = 1.054255712324226
$block1[6,0] = 1.02
$block1[6,1] = 1.038226637493707
$block1[6,2] = 1.045153090569665
$block1[6,3] = 1.036891657966809
$block1[6,4] = 1.052600032649875
$block1[7,0] = 1.02
$block1[7,1] = 1.035480625196141
$block1[7,2] = 1.042614693307997
$block1[7,3] = 1.034543577365336
$block1[7,4] = 1.049685004413007
$block1[8,0] = 1.02
$block1[8,1] = 1.033648193795389
$block1[8,2] = 1.040922512613412
$block1[8,3] = 1.032978448120655
$block1[8,4] = 1.047743232950656
$block1[9,0] = 1.02
$block1[9,1] = 1.032854288646367
$block1[9,2] = 1.040189783576391
$block1[9,3] = 1.032300778999502
$block1[9,4] = 1.04690277760435
$block1[10,0] = 1.02
$block1[10,1] = 1.032559326979584
$block1[10,2] = 1.039917613353741
$block1[10,3] = 1.032049067711786
$block1[10,4] = 1.046590644847504
$block1[11,0] = 1.02
$block1[11,1] = 1.032622600460565
$block1[11,2] = 1.039975994906394
$block1[11,3] = 1.032103060420629
$block1[11,4] = 1.046657596162598
$block1[12,0] = 1.02
$block1[12,1] = 1.032829908463172
$block1[12,2] = 1.040167285950998
$block1[12,3] = 1.032279972366528
$block1[12,4] = 1.046876975609269
$block1[13,0] = 1.02
$block1[13,1] = 1.032957628428225
$block1[13,2] = 1.04028514636289
$block1[13,3] = 1.032388974309726
$block1[13,4] = 1.047012149063842
$block1[14,0] = 1.02
$block1[14,1] = 1.033700872982883
$block1[14,2] = 1.040971141223062
$block1[14,3] = 1.033023423579829
$block1[14,4] = 1.047799018293702
$block1[15,0] = 1.02
$block1[15,1] = 1.034166968486347
$block1[15,2] = 1.041401445654125
$block1[15,3] = 1.033421407332932
$block1[15,4] = 1.048292691412365
$block1[16,0] = 1.02
$block1[16,1] = 1.034438790821059
$block1[16,2] = 1.041652434700495
$block1[16,3] = 1.033653548670026
$block1[16,4] = 1.048580675965756
$block1[17,0] = 1.02
$block1[17,1] = 1.03453146793947
$block1[17,2] = 1.041738015481522
$block1[17,3] = 1.033732703604221
$block1[17,4] = 1.048678876976599
$block1[18,0] = 1.02
$block1[18,1] = 1.034116965345062
$block1[18,2] = 1.041355278098465
$block1[18,3] = 1.033378707022878
$block1[18,4] = 1.048239721522081
$block1[19,0] = 1.02
$block1[19,1] = 1.032768863372763
$block1[19,2] = 1.040110955574286
$block1[19,3] = 1.032227876073586
$block1[19,4] = 1.046812372471862
$block1[20,0] = 1.02
$block1[20,1] = 1.031920852329009
$block1[20,2] = 1.03932858810936
$block1[20,3] = 1.031504332479761
$block1[20,4] = 1.045915228500459
$block1[21,0] = 1.02
$block1[21,1] = 1.032370438265658
$block1[21,2] = 1.039743337589767
$block1[21,3] = 1.031887894214524
$block1[21,4] = 1.046390794949938
$block1[22,0] = 1.02
$block1[22,1] = 1.034139559746031
$block1[22,2] = 1.041376139230476
$block1[22,3] = 1.033398001442285
$block1[22,4] = 1.048263656232391
$block1[23,0] = 1.02
$block1[23,1] = 1.03619083664287
$block1[23,2] = 1.043270909988038
$block1[23,3] = 1.035150563348635
$block1[23,4] = 1.050438323565232
$ws.Range("B2:F25").Value = $block1

$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 1.033188014343969
$block2[0,1] = 1.04292641433427
$block2[0,2] = 1.047552306555255
$block2[0,3] = 1.039342003069077
$block2[0,4] = 1.054923429018759
$block2[0,5] = 1.018169716802007
$block2[1,0] = 1.033356328596605
$block2[1,1] = 1.043756085295388
$block2[1,2] = 1.048460273666902
$block2[1,3] = 1.040166287665935
$block2[1,4] = 1.055995996225231
$block2[1,5] = 1.018450693363166
$block2[2,0] = 1.033463563108984
$block2[2,1] = 1.044292116754568
$block2[2,2] = 1.049047235671842
$block2[2,3] = 1.040699114538661
$block2[2,4] = 1.056689795518457
$block2[2,5] = 1.018632078910901
$block2[3,0] = 1.033508243297921
$block2[3,1] = 1.044517268656852
$block2[3,2] = 1.04929386294869
$block2[3,3] = 1.040922986324999
$block2[3,4] = 1.0569814164592
$block2[3,5] = 1.01870823182726
$block2[4,0] = 1.033515721773185
$block2[4,1] = 1.044555061199525
$block2[4,2] = 1.049335265066869
$block2[4,3] = 1.040960567861395
$block2[4,4] = 1.057030377867418
$block2[4,5] = 1.018721012286682
$block2[5,0] = 1.033464161704236
$block2[5,1] = 1.044295126013305
$block2[5,2] = 1.049050531633321
$block2[5,3] = 1.040702106428949
$block2[5,4] = 1.056693692375678
$block2[5,5] = 1.018633096868748
$block2[6,0] = 1.033245244006553
$block2[6,1] = 1.043206976057993
$block2[6,2] = 1.047859273451469
$block2[6,3] = 1.039620686241992
$block2[6,4] = 1.055285955111277
$block2[6,5] = 1.018264762279821
$block2[7,0] = 1.032846643681279
$block2[7,1] = 1.041283189020881
$block2[7,2] = 1.045755843023374
$block2[7,3] = 1.037710913512063
$block2[7,4] = 1.052803576818345
$block2[7,5] = 1.017612445066407
$block2[8,0] = 1.032572272491704
$block2[8,1] = 1.039996352997747
$block2[8,2] = 1.044350615334529
$block2[8,3] = 1.03643487719462
$block2[8,4] = 1.051147392428334
$block2[8,5] = 1.01717535830186
$block2[9,0] = 1.032451416412094
$block2[9,1] = 1.039438102759774
$block2[9,2] = 1.043741424085633
$block2[9,3] = 1.035881650575288
$block2[9,4] = 1.05042992781926
$block2[9,5] = 1.016985567675315
$block2[10,0] = 1.032406216658153
$block2[10,1] = 1.03923058599604
$block2[10,2] = 1.043515034041115
$block2[10,3] = 1.035676052131309
$block2[10,4] = 1.050163378884506
$block2[10,5] = 1.016914991151489
$block2[11,0] = 1.032415926118978
$block2[11,1] = 1.039275106205954
$block2[11,2] = 1.043563600441703
$block2[11,3] = 1.035720158484938
$block2[11,4] = 1.050220556854165
$block2[11,5] = 1.016930133673482
$block2[12,0] = 1.032447686479293
$block2[12,1] = 1.039420952585094
$block2[12,2] = 1.043722712851315
$block2[12,3] = 1.035864657901253
$block2[12,4] = 1.050407895843818
$block2[12,5] = 1.016979735428375
$block2[13,0] = 1.032467214210091
$block2[13,1] = 1.039510792409122
$block2[13,2] = 1.043820732735397
$block2[13,3] = 1.03595367475497
$block2[13,4] = 1.050523314818477
$block2[13,5] = 1.017010286117769
$block2[14,0] = 1.032580250070684
$block2[14,1] = 1.040033380171775
$block2[14,2] = 1.044391030099782
$block2[14,3] = 1.036471578297108
$block2[14,4] = 1.051195001183192
$block2[14,5] = 1.017187942903301
$block2[15,0] = 1.032650605079952
$block2[15,1] = 1.04036090587608
$block2[15,2] = 1.044748569277395
$block2[15,3] = 1.036796258794864
$block2[15,4] = 1.051616243872769
$block2[15,5] = 1.017299240402533
$block2[16,0] = 1.032691444036654
$block2[16,1] = 1.040551845650865
$block2[16,2] = 1.044957046610052
$block2[16,3] = 1.036985572505644
$block2[16,4] = 1.051861916282988
$block2[16,5] = 1.017364107316294
$block2[17,0] = 1.032705335499113
$block2[17,1] = 1.040616934183108
$block2[17,2] = 1.045028120292598
$block2[17,3] = 1.037050112269196
$block2[17,4] = 1.051945678963343
$block2[17,5] = 1.017386216611091
$block2[18,0] = 1.032643077127878
$block2[18,1] = 1.040325775865292
$block2[18,2] = 1.044710215877902
$block2[18,3] = 1.036761430581951
$block2[18,4] = 1.051571051771921
$block2[18,5] = 1.01728730451561
$block2[19,0] = 1.032438342360469
$block2[19,1] = 1.0393780088269
$block2[19,2] = 1.043675861240366
$block2[19,3] = 1.035822109344433
$block2[19,4] = 1.050352730591931
$block2[19,5] = 1.016965131154362
$block2[20,0] = 1.032307832932651
$block2[20,1] = 1.038781196641598
$block2[20,2] = 1.043024888427553
$block2[20,3] = 1.035230910174847
$block2[20,4] = 1.049586429923147
$block2[20,5] = 1.016762105909664
$block2[21,0] = 1.032377187673392
$block2[21,1] = 1.03909766507774
$block2[21,2] = 1.043370041889576
$block2[21,3] = 1.035544374292578
$block2[21,4] = 1.049992688817081
$block2[21,5] = 1.016869777343463
$block2[22,0] = 1.032646479297007
$block2[22,1] = 1.040341649915391
$block2[22,2] = 1.044727546342188
$block2[22,3] = 1.036777168159093
$block2[22,4] = 1.051591472233883
$block2[22,5] = 1.017292697986865
$block2[23,0] = 1.032951213473557
$block2[23,1] = 1.041781289794736
$block2[23,2] = 1.046300143952108
$block2[23,3] = 1.03820513502674
$block2[23,4] = 1.053445548387657
$block2[23,5] = 1.017781472973839
$ws.Range("I2:N25").Value = $block2

Write-Host "applied vm_pu updates"
